$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 908, shifting existing rows 908:1005 down to 910:1007
$ws.Rows("908:909").Insert()

# Fill in the constant columns (same for every data row in this sheet)
$mercadoId = 5
$mercado = "Macroferia Regional de Talca"
$region = "Maule"
$codreg = 7
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102005
$categoria = "Naranja"

# New row 908: Naranja, Lane Late, Primera
$ws.Range("A908").Value = $mercadoId
$ws.Range("B908").Value = $mercado
$ws.Range("C908").Value = $region
$ws.Range("D908").Value = 45194
$ws.Range("E908").Value = $codreg
$ws.Range("F908").Value = $tipo
$ws.Range("G908").Value = $productoId
$ws.Range("H908").Value = $producto
$ws.Range("I908").Value = $categoriaId
$ws.Range("J908").Value = $categoria
$ws.Range("K908").Value = "Lane Late"
$ws.Range("L908").Value = "Primera"
$ws.Range("M908").Value = 250
$ws.Range("N908").Value = 9000
$ws.Range("O908").Value = 9000
$ws.Range("P908").Value = 9000
$ws.Range("Q908").Value = "$/bandeja 15 kilos granel"
$ws.Range("R908").Value = "Región de O'Higgins"
$ws.Range("S908").Value = 600
$ws.Range("T908").Value = 15

# New row 909: Naranja, Navel Late, Primera
$ws.Range("A909").Value = $mercadoId
$ws.Range("B909").Value = $mercado
$ws.Range("C909").Value = $region
$ws.Range("D909").Value = 45194
$ws.Range("E909").Value = $codreg
$ws.Range("F909").Value = $tipo
$ws.Range("G909").Value = $productoId
$ws.Range("H909").Value = $producto
$ws.Range("I909").Value = $categoriaId
$ws.Range("J909").Value = $categoria
$ws.Range("K909").Value = "Navel Late"
$ws.Range("L909").Value = "Primera"
$ws.Range("M909").Value = 300
$ws.Range("N909").Value = 9000
$ws.Range("O909").Value = 9000
$ws.Range("P909").Value = 9000
$ws.Range("Q909").Value = "$/bandeja 15 kilos granel"
$ws.Range("R909").Value = "Región de O'Higgins"
$ws.Range("S909").Value = 600
$ws.Range("T909").Value = 15
